$d = $word.ActiveDocument

# Locate the last paragraph (ends with the run containing "呀") and
# append a new run "！" right after it, inheriting the same run
# formatting (宋体, size 28) used throughout the document.
$paras = $d.Paragraphs
$lastPara = $paras.Last

$r = $lastPara.Range
# Collapse to the end of the paragraph's text (before the paragraph mark)
$insertRange = $d.Range($r.End - 1, $r.End - 1)
$insertRange.InsertAfter("！")

# Apply matching character formatting to the newly inserted text.
$newRange = $d.Range($r.End - 1, $r.End)
$newRange.Font.Name = "宋体"
$newRange.Font.NameFarEast = "宋体"
$newRange.Font.Size = 14
